# Updated symbol list on Sat Dec 31 09:36:12 UTC 2022 with GitHub Actions
# Refresh the cryptocurrency price/volume snapshot on Sheet1.
# Numeric-looking "Price" values live in column D as text (e.g. "245.80"),
# so they are written with a leading apostrophe to keep them stored as
# text rather than being coerced into numeric cells. The "Volume(1h)"
# labels in column E are plain text already and need no such treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'245.47"
$ws.Range("D3").Value  = "'26.10"
$ws.Range("D4").Value  = "'5.120"
$ws.Range("D5").Value  = "'0.05589"
$ws.Range("D6").Value  = "'6.494"
$ws.Range("D7").Value  = "'3.027"
$ws.Range("D8").Value  = "'0.8168"
$ws.Range("D9").Value  = "'0.8476"
$ws.Range("D10").Value = "'0.1337"
$ws.Range("D12").Value = "'0.02852"
$ws.Range("D13").Value = "'0.09396"
$ws.Range("D14").Value = "'0.001518"

$ws.Range("D15").Value = "'0.0006010"
$ws.Range("E15").Value = "14OneONEWorstin24h"

$ws.Range("D16").Value = "'0.006230"
$ws.Range("D17").Value = "'3.551"
$ws.Range("D20").Value = "'0.06946"
$ws.Range("D22").Value = "'3.740"
$ws.Range("D23").Value = "'0.04688"
$ws.Range("D25").Value = "'0.001249"
$ws.Range("D26").Value = "'0.004603"
$ws.Range("D27").Value = "'0.00009603"
$ws.Range("D40").Value = "'0.03653"

$ws.Range("D41").Value = "'0.003398"

$ws.Range("D42").Value = "'0.1350"
$ws.Range("E42").Value = "41BKEXTokenBKKBestin24h"

$ws.Range("D44").Value = "'0.007889"
$ws.Range("D45").Value = "'0.00005313"

$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
